$d = $word.ActiveDocument

# The target paragraph currently holds "<id>p145v_2</id>" split across three
# runs: "<id>" (Courier New / 7f6000 / 18pt), "p145v_2" (default formatting),
# and "</id>" (Courier New / 7f6000 / 18pt again). The edit collapses these
# into a single run - using the first run's formatting - that contains the
# whole string "<id>p145v_2</id>".

$search = $d.Content
$found = $search.Find.Execute(
    "<id>p145v_2</id>", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
if (-not $found) {
    throw "Could not find target text '<id>p145v_2</id>'"
}

$matchStart = $search.Start
$matchEnd = $search.End

# Isolate the leading "<id>" run - we keep it untouched so it keeps its
# original run properties (and serialization, e.g. xml:space="preserve").
$openTag = "<id>"
$firstRun = $d.Range($matchStart, $matchStart + $openTag.Length)
if ($firstRun.Text -ne $openTag) {
    throw "Unexpected content at start of match: $($firstRun.Text)"
}

# Grab the remaining text ("p145v_2</id>", i.e. runs 2 and 3), delete those
# runs, then re-insert the text right after the untouched first run so Word
# merges it into that run instead of minting fresh, differently-formatted
# runs.
$remainder = $d.Range($matchStart + $openTag.Length, $matchEnd)
$remainderText = $remainder.Text
$remainder.Text = ""

$firstRun = $d.Range($matchStart, $matchStart + $openTag.Length)
$firstRun.InsertAfter($remainderText)
